# Add a new "Save" column (H) to the sheet, mirroring the header style
# used by the existing stat columns, and populate it with the per-row
# save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same text/value style as the other headers (e.g. G1 "sum")
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row "Save" values for rows 2-16 (H2:H16)
$saveValues = @(0, 1, 1, 1, 1, 0, 0, 1, 0, 0, 0, 0, 1, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
